$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row (B6) before updating B5 so that the shared-string table
# ends up in the same order as the target file (new string for B6 first,
# then the replacement string for B5).
$ws.Range("B6").Value = "In allen Tabellen Zeilen rausschmeißen, zu denen es keine passende Diagnosis- und Patient-Data gibtr"

# Replace the text in B5 ("Value transformation mapping aus externer Datei"
# -> "Duplikate in allen Tabellen").
$ws.Range("B5").Value = "Duplikate in allen Tabellen"

# Widen column B.
$ws.Columns.Item(2).ColumnWidth = 97

# Match the row height formatting used by the other data rows.
$ws.Rows.Item(6).RowHeight = 25.2

# Move/extend the active selection to D16, matching the saved view state.
$ws.Range("D16").Select()
